# Auto-generated edit script: updates crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell 2 4 "42.733.01"
Set-TextCell 2 5 "  -1.79%  "
Set-TextCell 3 4 "2.332.89"
Set-TextCell 3 5 "  +0.51%  "
Set-TextCell 4 4 "0.999"
Set-TextCell 4 5 "  -0.07%  "
Set-TextCell 5 4 "306.16"
Set-TextCell 5 5 "  -1.93%  "
Set-TextCell 6 4 "100.29"
Set-TextCell 6 5 "  -2.23%  "
Set-TextCell 7 5 "  -5.38%  "
Set-TextCell 8 5 "  +0.00%  "
Set-TextCell 9 4 "0.509"
Set-TextCell 9 5 "  -4.93%  "
Set-TextCell 10 4 "34.91"
Set-TextCell 10 5 "  -3.25%  "
Set-TextCell 11 4 "52.10"
Set-TextCell 11 5 "  +0.27%  "
Set-TextCell 12 5 "  -2.30%  "
Set-TextCell 13 5 "  -0.87%  "
Set-TextCell 14 5 "  -3.70%  "
Set-TextCell 15 4 "15.73"
Set-TextCell 15 5 "  +4.66%  "
Set-TextCell 16 4 "2.278.08"
Set-TextCell 16 5 "  -1.39%  "
Set-TextCell 17 4 "0.795"
Set-TextCell 17 5 "  -2.48%  "
Set-TextCell 18 4 "42.666.27"
Set-TextCell 18 5 "  -1.72%  "
Set-TextCell 19 4 "6.26"
Set-TextCell 19 5 "  +1.57%  "
Set-TextCell 20 4 "0.0₃0904"
Set-TextCell 20 5 "  -2.53%  "
Set-TextCell 21 4 "11.59"
Set-TextCell 21 5 "  -7.44%  "
Set-TextCell 22 4 "67.65"
Set-TextCell 22 5 "  -1.14%  "
Set-TextCell 23 4 "236.29"
Set-TextCell 23 5 "  -2.61%  "
Set-TextCell 24 4 "1.99"
Set-TextCell 24 5 "  -3.14%  "
Set-TextCell 25 4 "2.56"
Set-TextCell 25 5 "  -2.82%  "
Set-TextCell 26 5 "  -0.07%  "
Set-TextCell 27 4 "24.83"
Set-TextCell 27 5 "  +0.08%  "
Set-TextCell 28 4 "2.19"
Set-TextCell 28 5 "  +3.61%  "
Set-TextCell 29 4 "34.83"
Set-TextCell 29 5 "  -7.64%  "
Set-TextCell 30 4 "9.34"
Set-TextCell 30 5 "  -3.43%  "
Set-TextCell 31 4 "159.08"
Set-TextCell 31 5 "  -4.81%  "
Set-TextCell 32 4 "0.998"
Set-TextCell 32 5 "  -0.15%  "
Set-TextCell 33 5 "  -4.38%  "
Set-TextCell 34 2 "WEMIXToken"
Set-TextCell 34 3 "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell 34 4 "2.45"
Set-TextCell 34 5 "  -3.07%  "
Set-TextCell 35 2 "Celestia"
Set-TextCell 35 3 "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextCell 35 4 "17.37"
Set-TextCell 35 5 "  -2.44%  "
Set-TextCell 36 4 "0.0726"
Set-TextCell 36 5 "  -3.08%  "
Set-TextCell 37 4 "4.56"
Set-TextCell 37 5 "  +5.32%  "
Set-TextCell 38 4 "2.96"
Set-TextCell 38 5 "  -4.93%  "
Set-TextCell 39 5 "  -0.84%  "
Set-TextCell 40 5 "  -4.31%  "
Set-TextCell 41 5 "  -3.69%  "
Set-TextCell 42 4 "2.35"
Set-TextCell 42 5 "  +1.46%  "
Set-TextCell 43 4 "2.013.10"
Set-TextCell 43 5 "  +1.42%  "
Set-TextCell 44 4 "0.0284"
Set-TextCell 44 5 "  -2.36%  "
Set-TextCell 45 4 "18.78"
Set-TextCell 45 5 "  -5.27%  "
Set-TextCell 46 4 "10.28"
Set-TextCell 46 5 "  +3.87%  "
Set-TextCell 47 4 "2.92"
Set-TextCell 47 5 "  -3.25%  "
Set-TextCell 48 4 "55.62"
Set-TextCell 48 5 "  -0.71%  "
Set-TextCell 49 4 "2.92"
Set-TextCell 49 5 "  -0.68%  "
Set-TextCell 50 4 "2.558.71"
Set-TextCell 50 5 "  +0.50%  "
Set-TextCell 51 4 "4.63"
Set-TextCell 51 5 "  +2.03%  "
